$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the "total" row (currently row 4), shifting it to row 5
$ws.Rows.Item(4).Insert()

# New "ops" row
$ws.Range("A4").Value = "ops"
$ws.Range("B4").Value = 20
$ws.Range("D4").Value = "ops"
$ws.Range("E4").Value = 20

# Fix up the total formulas to include the new row
$ws.Range("B5").Formula = "=SUM(B2:B4)"
$ws.Range("E5").Formula = "=SUM(E2:E4)"

# Update selection to match target state
$ws.Range("F3").Select()
